# Added New Mac-Address and Document Types
# Append a new data row (row 33) to the master-reg_center_user_h sheet,
# mirroring the existing rows' pattern, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 33
$ws.Cells.Item($row, 1).Value = 10002        # regcntr_id
$ws.Cells.Item($row, 2).Value = 110032       # usr_id
$ws.Cells.Item($row, 3).Value = "eng"        # lang_code
$ws.Cells.Item($row, 4).Value = $true        # is_active
$ws.Cells.Item($row, 5).Value = "superadmin" # cr_by
$ws.Cells.Item($row, 6).Value = "now()"      # cr_dtimes
$ws.Cells.Item($row, 7).Value = "now()"      # eff_dtimes

# Move the selected/active cell as recorded in the saved view state.
$null = $ws.Range("D26").Select()
